# Apply updated cryptocurrency data (prices, 1h volume %, and row reorder for Celestia/Algorand/LidoDAOToken)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.612.55"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").Value = "2.290.74"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'114.60"
$ws.Range("E5").Value = "  +19.05%  "
$ws.Range("D6").Value = "'269.11"
$ws.Range("E6").Value = "  +1.09%  "
$ws.Range("D7").Value = "'0.623"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("D9").Value = "'0.615"
$ws.Range("E9").Value = "  +1.27%  "
$ws.Range("D10").Value = "'47.97"
$ws.Range("E10").Value = "  +5.14%  "
$ws.Range("D11").Value = "'0.0939"
$ws.Range("E11").Value = "  +0.64%  "
$ws.Range("E12").Value = "  +10.68%  "
$ws.Range("D13").Value = "'0.106"
$ws.Range("E13").Value = "  +0.88%  "
$ws.Range("D14").Value = "'15.58"
$ws.Range("E14").Value = "  +2.99%  "
$ws.Range("D15").Value = "2.635.78"
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("D16").Value = "'0.848"
$ws.Range("E16").Value = "  +0.45%  "
$ws.Range("D17").Value = "2.293.18"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("D18").Value = "43.588.87"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("E19").Value = "  +2.44%  "
$ws.Range("D20").Value = "'6.56"
$ws.Range("E20").Value = "  +5.79%  "
$ws.Range("D21").Value = "'72.53"
$ws.Range("E21").Value = "  +0.96%  "
$ws.Range("E22").Value = "  +2.60%  "
$ws.Range("D23").Value = "'233.38"
$ws.Range("E23").Value = "  +0.52%  "
$ws.Range("D24").Value = "'9.52"
$ws.Range("E24").Value = "  +4.61%  "
$ws.Range("D25").Value = "'2.84"
$ws.Range("E25").Value = "  +14.36%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").Value = "'11.55"
$ws.Range("E27").Value = "  +3.94%  "
$ws.Range("D28").Value = "'42.01"
$ws.Range("E28").Value = "  +5.29%  "
$ws.Range("D29").Value = "'3.40"
$ws.Range("E29").Value = "  -1.84%  "
$ws.Range("E30").Value = "  +0.21%  "
$ws.Range("D31").Value = "'176.66"
$ws.Range("E31").Value = "  +0.72%  "
$ws.Range("D32").Value = "'21.61"
$ws.Range("E32").Value = "  -0.78%  "
$ws.Range("D33").Value = "'0.0930"
$ws.Range("E33").Value = "  +5.35%  "
$ws.Range("E34").Value = "  +3.70%  "
$ws.Range("E35").Value = "  +0.85%  "
$ws.Range("E36").Value = "  +10.02%  "
$ws.Range("E37").Value = "  +0.58%  "
$ws.Range("E38").Value = "  +1.40%  "
$ws.Range("E39").Value = "  +13.57%  "
$ws.Range("D40").Value = "'74.04"
$ws.Range("E40").Value = "  +15.28%  "
$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D41").Value = "'13.82"
$ws.Range("E41").Value = "  +12.82%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "'0.243"
$ws.Range("E42").Value = "  +3.75%  "
$ws.Range("B43").Value = "LidoDAOToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D43").Value = "'2.39"
$ws.Range("E43").Value = "  +4.01%  "
$ws.Range("D44").Value = "'1.43"
$ws.Range("E44").Value = "  +6.71%  "
$ws.Range("E45").Value = "  +14.63%  "
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").Value = "'8.76"
$ws.Range("E47").Value = "  -0.22%  "
$ws.Range("D48").Value = "'102.53"
$ws.Range("E48").Value = "  +5.38%  "
$ws.Range("E49").Value = "  -1.49%  "
$ws.Range("D50").Value = "'1.24"
$ws.Range("E50").Value = "  +4.55%  "
$ws.Range("E51").Value = "  +4.39%  "
